# Balance : Location selection dynamic
# Replace the numeric "To/From" location codes (101-111) in column E with
# descriptive location names, and extend the used range with a few extra
# blank (but formatted) rows so the location column keeps its formatting
# further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New location names for column E, row by row (rows 2-16)
$locations = @{
    2  = "WA Nursing Home"
    3  = "External facility"
    4  = "Facility 3"
    5  = "Wing A"
    6  = "WA Nursing Home"
    7  = "sub loc 2"
    8  = "Integration Level 1"
    9  = "Integration level 1"
    10 = "Room 5"
    11 = "WA Nursing Home"
    12 = "External facility"
    13 = "Facility 3"
    14 = "Wing A"
    15 = "WA Nursing Home"
    16 = "sub loc 1"
}

foreach ($row in 2..16) {
    $ws.Cells.Item($row, 5).Value = $locations[$row]
}

# Give the location column a left/bottom alignment (previously it relied on
# the default alignment) and extend the same formatting a few rows further
# down (17-19) so new entries keep the same look even before data is typed.
$ws.Range("E2:E19").HorizontalAlignment = -4131
$ws.Range("E2:E19").VerticalAlignment = -4107

$ws.Range("E17").Select()
